# Distribution channel add for all module excel upload
#
# Inserts a new "Distribution channel code" column before the existing
# "Budget" column (column I), pushing Budget to column J, and fills the
# new column with the distribution channel codes for the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at I; this shifts the old column I ("Budget") to J
# and carries over formatting, exactly like a manual "Insert Column".
$null = $ws.Columns("I:I").Insert()

# New header + data for the inserted "Distribution channel code" column.
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Give the new column a sensible custom width so the header text fits.
$ws.Columns("I:I").ColumnWidth = 21.6

# Restore/update the active selection the way Excel would after an
# insert performed with the cursor further to the right of the sheet.
$null = $ws.Range("M13").Select()
